$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) columns.
# NumberFormat is forced to text ("@") before writing D-column values so
# strings such as "27.45", "1.80" or "0.0000100" are not coerced into
# numbers (which would silently drop significant trailing zeros / switch
# to scientific notation). Style is reset to "Normal" immediately after
# so the cell keeps its original (unstyled) appearance.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.802.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.598.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.88%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.551"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.594.75"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.121"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.87%  "
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.345"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.077.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000178"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.710.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.605.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "364.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.36%  "
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "67.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.739.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "583.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000100"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.122"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "155.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.364"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "153.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0290"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.612"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.19%  "
